$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly to fit the new, longer descriptions (~29.7 -> ~33.1 chars)
$ws.Columns.Item(3).ColumnWidth = 32.43

# New row 13: a "ユーザー編集画面" (user-edit screen) record, transitioning from
# the existing "ユーザー一覧画面" (user-list) screen added previously.
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "ユーザー編集画面"
$ws.Cells.Item(13, 3).Value = "管理者がユーザー情報を編集する画面"
$ws.Cells.Item(13, 4).Value = "ユーザー一覧画面"

# Copy the ID column's number formatting/font down onto the new row
# (xlPasteFormats = -4122; leaves the value we just entered untouched)
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection to the newly-entered cell, as Excel would leave it
$ws.Range("D13").Select()
